$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$status = "Ready for handoff"
$overviewDate = "2016-09-03 16:55:20"
$zhcnHandoffDate = "2016-09-03 16:55:15"
$dedeHandoffDate = "2016-09-03 16:55:20"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/149bcc218bd00c6efe143eda8dcc0568e4b1bcb9/e2e/d86bf385-42c1-476a-8eb2-d5c78d48af64.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c0c92dfd8d0786b36ce5a1a85490efff0b2642ad/e2e/d86bf385-42c1-476a-8eb2-d5c78d48af64.md."

# Overview sheet: row 3 corresponds to d86bf385-...md
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status
$overview.Range("G3").Value = $overviewDate

# zh-cn sheet
$zhcn.Range("C2").Value = $status
$zhcn.Range("C3").Value = $status
$zhcn.Range("H3").Value = $zhcnHandoffDate
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# de-de sheet
$dede.Range("C2").Value = $status
$dede.Range("C3").Value = $status
$dede.Range("H3").Value = $dedeHandoffDate
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.17
